$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.043.94"
$ws.Range("E2").Value = "  -2.76%  "
$ws.Range("D3").Value = "2.339.69"
$ws.Range("E3").Value = "  -3.71%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.58"
$ws.Range("E5").Value = "  -2.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "85.50"
$ws.Range("E6").Value = "  -3.68%  "
$ws.Range("E7").Value = "  -2.14%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.485"
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0813"
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.03"
$ws.Range("E11").Value = "  -6.44%  "
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").Value = "2.686.05"
$ws.Range("E13").Value = "  -4.22%  "
$ws.Range("E14").Value = "  -4.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.79"
$ws.Range("E15").Value = "  -5.05%  "
$ws.Range("D16").Value = "2.354.86"
$ws.Range("E16").Value = "  -2.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.757"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").Value = "40.038.35"
$ws.Range("E18").Value = "  -2.66%  "
$ws.Range("E19").Value = "  -2.07%  "
$ws.Range("E20").Value = "  -1.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.94"
$ws.Range("E21").Value = "  -5.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.66"
$ws.Range("E22").Value = "  -3.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.38"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  -5.33%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.35"
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.31"
$ws.Range("E29").Value = "  -2.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.88"
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.20"
$ws.Range("E31").Value = "  -2.52%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.11"
$ws.Range("E33").Value = "  -2.74%  "
$ws.Range("E34").Value = "  -3.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0721"
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.81"
$ws.Range("E37").Value = "  -3.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0991"
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.70"
$ws.Range("E39").Value = "  -5.19%  "
$ws.Range("E40").Value = "  -2.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.86"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "1.958.88"
$ws.Range("E42").Value = "  -1.18%  "
$ws.Range("E43").Value = "  -4.33%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0264"
$ws.Range("E44").Value = "  -3.98%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.54"
$ws.Range("E45").Value = "  -4.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.41"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.70"
$ws.Range("E47").Value = "  -5.66%  "
$ws.Range("D48").Value = "2.561.45"
$ws.Range("E48").Value = "  -4.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.71"
$ws.Range("E49").Value = "  -2.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.70"
$ws.Range("E50").Value = "  -3.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.69"
$ws.Range("E51").Value = "  -2.27%  "
